$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark. In the original
#     document it sits right before the "ASKED HARTMUT AND ABE TO
#     CLARIFY" paragraph; it needs to move to wrap the "XXX"
#     placeholder inside the breakdown-voltage sentence instead. ---
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# --- Step 2: find the unique sentence containing the "XXX" placeholder
#     so we don't accidentally match one of the other "XXX" markers
#     elsewhere in the document. ---
$sentence = $d.Content.Duplicate
$sentence.Find.Execute("The breakdown voltage of HPK 50D sensor is about 650 V, and that of CNM is XXX.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# --- Step 3: within that sentence, grab just the "XXX" run and make it
#     bold (the rest of the sentence keeps its original formatting). ---
$xxxRng = $sentence.Duplicate
$xxxRng.Find.Execute("XXX", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xxxRng.Bold = 1

# --- Step 4: re-create the "_GoBack" bookmark around the now-bold
#     "XXX" run. ---
$d.Bookmarks.Add("_GoBack", $xxxRng) | Out-Null
